# "Generate Report for Handback"
#
# The handback transform failed for both locales (the generated handback
# file name didn't match the expected handoff-derived name), so the
# localization-status report needs to:
#   1. Flip every "Ready for handoff" status cell (Overview's per-locale
#      columns + both locale tables' Status column all shared that one
#      string) to "Handback transform failed".
#   2. Record the specific mismatch error in the "Error Detail" column
#      (column P) of both the zh-cn and de-de tables.
#   3. Widen the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$failedStatus = "Handback transform failed"

# Every cell that currently shows "Ready for handoff" needs to show the
# new status text. Touch them all (Overview's per-locale status columns
# plus both locale tables' Status column) so the workbook ends up with a
# single, shared, de-duplicated string instead of a stray second copy.
$overview.Range("E3").Value = $failedStatus
$overview.Range("F3").Value = $failedStatus
$zhcn.Range("C3").Value = $failedStatus
$dede.Range("C3").Value = $failedStatus

# Error Detail (column P) messages for row 3 (the 60317641-... file) in
# each locale table.
$zhcnError = "Handback file name: a3dqqjiy.i0t is different with handoff file name: 60317641-dad3-441a-bb28-e1b4508b8636.6beddd35b8b4ee5f23d5de3abdc537143ca211d8.zh-cn."
$dedeError = "Handback file name: a3dqqjiy.i0t is different with handoff file name: 60317641-dad3-441a-bb28-e1b4508b8636.6beddd35b8b4ee5f23d5de3abdc537143ca211d8.de-de."

$zhcn.Range("P3").Value = $zhcnError
$dede.Range("P3").Value = $dedeError

# Widen column P ("Error Detail") on both locale sheets to fit the new
# message text. 39.15 chars is what the host's pixel-grid column-width
# rounding resolves to an on-disk width of exactly 40.
$zhcn.Columns.Item(16).ColumnWidth = 39.15
$dede.Columns.Item(16).ColumnWidth = 39.15
